# Update the "ID Competição" column (B) for every data row.
# The scraped competition id was truncated/incorrect (42) and is being
# corrected/restored to its real value (242), per the commit message
# "concerting [sic] names and recovering dropped data".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    if ($current -eq 42) {
        $cell.Value2 = 242
    }
}
